$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp (shared string table last entry)
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 14:50"

# Row 8: Alemania
$ws.Range("A8").Value = "Alemania"
$ws.Range("B8").Value = 63079
$ws.Range("C8").Value = 644
$ws.Range("D8").Value = 9211
$ws.Range("E8").Value = 53323
$ws.Range("F8").Value = 1979
$ws.Range("G8").Value = 4
$ws.Range("H8").Value = 545

# Row 14: Paises Bajos
$ws.Range("A14").Value = "Paises Bajos"
$ws.Range("B14").Value = 11750
$ws.Range("C14").Value = 884
$ws.Range("D14").Value = 250
$ws.Range("E14").Value = 10636
$ws.Range("F14").Value = 1053
$ws.Range("G14").Value = 93
$ws.Range("H14").Value = 864

# Row 25: Chequia
$ws.Range("A25").Value = "Chequia"
$ws.Range("B25").Value = 2878
$ws.Range("C25").Value = 61
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 2850
$ws.Range("F25").Value = 52
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 17

# Row 48: Singapur
$ws.Range("A48").Value = "Singapur"
$ws.Range("B48").Value = 879
$ws.Range("C48").Value = 35
$ws.Range("D48").Value = 228
$ws.Range("E48").Value = 648
$ws.Range("F48").Value = 19
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 3

# Row 49: Republica Dominicana
$ws.Range("A49").Value = "Republica Dominicana"
$ws.Range("B49").Value = 859
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 3
$ws.Range("E49").Value = 817
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 39

# Row 50: Peru
$ws.Range("A50").Value = "Peru"
$ws.Range("B50").Value = 852
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 16
$ws.Range("E50").Value = 818
$ws.Range("F50").Value = 40
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 18

# Row 51: Argentina
$ws.Range("A51").Value = "Argentina"
$ws.Range("B51").Value = 820
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 91
$ws.Range("E51").Value = 707
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 2
$ws.Range("H51").Value = 22

# Row 52: Croacia
$ws.Range("A52").Value = "Croacia"
$ws.Range("B52").Value = 790
$ws.Range("C52").Value = 77
$ws.Range("D52").Value = 67
$ws.Range("E52").Value = 717
$ws.Range("F52").Value = 27
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 6

# Row 54: Serbia
$ws.Range("A54").Value = "Serbia"
$ws.Range("B54").Value = 741
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 42
$ws.Range("E54").Value = 684
$ws.Range("F54").Value = 25
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 15

# Row 91: Vietnam
$ws.Range("A91").Value = "Vietnam"
$ws.Range("B91").Value = 203
$ws.Range("C91").Value = 9
$ws.Range("D91").Value = 55
$ws.Range("E91").Value = 148
$ws.Range("F91").Value = 3
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 0

# Row 126: Kenia
$ws.Range("A126").Value = "Kenia"
$ws.Range("B126").Value = 50
$ws.Range("C126").Value = 8
$ws.Range("D126").Value = 1
$ws.Range("E126").Value = 48
$ws.Range("F126").Value = 2
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 1

# Row 127: Banglades
$ws.Range("A127").Value = "Banglades"
$ws.Range("B127").Value = 49
$ws.Range("C127").Value = 1
$ws.Range("D127").Value = 19
$ws.Range("E127").Value = 25
$ws.Range("F127").Value = 1
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 5

# Row 128: Isla de Man
$ws.Range("A128").Value = "Isla de Man"
$ws.Range("B128").Value = 46
$ws.Range("C128").Value = 4
$ws.Range("D128").Value = 0
$ws.Range("E128").Value = 46
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 0

# Row 129: Monaco
$ws.Range("A129").Value = "Monaco"
$ws.Range("B129").Value = 46
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 1
$ws.Range("E129").Value = 44
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 1

# Row 130: Guayana Francesa
$ws.Range("A130").Value = "Guayana Francesa"
$ws.Range("B130").Value = 43
$ws.Range("C130").Value = 0
$ws.Range("D130").Value = 6
$ws.Range("E130").Value = 37
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 0

# Row 134: Guatemala
$ws.Range("A134").Value = "Guatemala"
$ws.Range("B134").Value = 36
$ws.Range("C134").Value = 2
$ws.Range("D134").Value = 10
$ws.Range("E134").Value = 25
$ws.Range("F134").Value = 1
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 1

# Row 135: Polinesia Francesa
$ws.Range("A135").Value = "Polinesia Francesa"
$ws.Range("B135").Value = 35
$ws.Range("C135").Value = 5
$ws.Range("D135").Value = 0
$ws.Range("E135").Value = 35
$ws.Range("F135").Value = 2
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 0

# Row 136: Jamaica
$ws.Range("A136").Value = "Jamaica"
$ws.Range("B136").Value = 34
$ws.Range("C136").Value = 2
$ws.Range("D136").Value = 2
$ws.Range("E136").Value = 31
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 1

# Row 138: Uganda
$ws.Range("A138").Value = "Uganda"
$ws.Range("B138").Value = 33
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 0
$ws.Range("E138").Value = 33
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 0

# Row 139: Guam
$ws.Range("A139").Value = "Guam"
$ws.Range("B139").Value = 32
$ws.Range("C139").Value = 0
$ws.Range("D139").Value = 0
$ws.Range("E139").Value = 31
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 1
